$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1676.8572
$ws.Range("J17").Value = 1795.5555
$ws.Range("L17").Value = 5386.666499999999
$ws.Range("N17").Value = -5722.666499999999
$ws.Range("H33").Value = 244.5
$ws.Range("I33").Value = 253.3077
$ws.Range("K33").Value = 253.3077
$ws.Range("M33").Value = -24.30770000000001
$ws.Range("H86").Value = 1201.1428
$ws.Range("I86").Value = 1201.1428
$ws.Range("K86").Value = 1201.1428
$ws.Range("M86").Value = -78.14280000000008
$ws.Range("H89").Value = 1201.1428
$ws.Range("I89").Value = 1201.1428
$ws.Range("K89").Value = 6005.714
$ws.Range("M89").Value = -389.7139999999999
$ws.Range("H112").Value = 1149.2354
$ws.Range("J112").Value = 1170.1613
$ws.Range("L112").Value = 3510.4839
$ws.Range("N112").Value = -5726.4839
$ws.Range("H131").Value = 3267.5
$ws.Range("I131").Value = 2523.3333
$ws.Range("J131").Value = 5500
$ws.Range("K131").Value = 7569.999899999999
$ws.Range("L131").Value = 16500
$ws.Range("M131").Value = -2529.999899999999
$ws.Range("N131").Value = -26580
$ws.Range("H132").Value = 3242.2363
$ws.Range("I132").Value = 1881.711
$ws.Range("J132").Value = 9364.6
$ws.Range("K132").Value = 5645.133
$ws.Range("L132").Value = 28093.8
$ws.Range("M132").Value = -3115.133
$ws.Range("N132").Value = -33153.8
$ws.Range("H137").Value = 73749090
$ws.Range("I137").Value = 500000740
$ws.Range("K137").Value = 1500002220
$ws.Range("M137").Value = -1499999670
$ws.Range("H138").Value = 2267.7083
$ws.Range("J138").Value = 2613.7646
$ws.Range("L138").Value = 7841.293799999999
$ws.Range("N138").Value = -18121.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3241.75
$ws.Range("I2").Value = 1890.4
$ws.Range("K2").Value = 1890.4
$ws.Range("M2").Value = -1777.4
$ws.Range("H32").Value = 2222950.2
$ws.Range("I32").Value = 2222950.2
$ws.Range("K32").Value = 2222950.2
$ws.Range("M32").Value = -2222663.2
$ws.Range("H61").Value = 2281.1667
$ws.Range("I61").Value = 1705.9565
$ws.Range("J61").Value = 4171.143
$ws.Range("K61").Value = 1705.9565
$ws.Range("L61").Value = 4171.143
$ws.Range("M61").Value = -1493.9565
$ws.Range("N61").Value = -4595.143
$ws.Range("H74").Value = 1988.6296
$ws.Range("I74").Value = 1175.9048
$ws.Range("J74").Value = 4833.1665
$ws.Range("K74").Value = 1175.9048
$ws.Range("L74").Value = 4833.1665
$ws.Range("M74").Value = -301.9048
$ws.Range("N74").Value = -6581.1665
$ws.Range("H77").Value = 1988.6296
$ws.Range("I77").Value = 1175.9048
$ws.Range("J77").Value = 4833.1665
$ws.Range("K77").Value = 5879.524
$ws.Range("L77").Value = 24165.8325
$ws.Range("M77").Value = -1511.524
$ws.Range("N77").Value = -32901.8325
$ws.Range("H97").Value = 1333.6471
$ws.Range("I97").Value = 1362.25
$ws.Range("K97").Value = 1362.25
$ws.Range("M97").Value = -866.25
$ws.Range("H110").Value = 1768.4375
$ws.Range("I110").Value = 1757.9166
$ws.Range("K110").Value = 1757.9166
$ws.Range("M110").Value = 287.0834
$ws.Range("H116").Value = 3241.75
$ws.Range("I116").Value = 1890.4
$ws.Range("K116").Value = 1890.4
$ws.Range("M116").Value = 403.5999999999999
$ws.Range("H122").Value = 1624.2174
$ws.Range("I122").Value = 1267.9
$ws.Range("K122").Value = 3803.7
$ws.Range("M122").Value = -1353.7
$ws.Range("H136").Value = 2281.1667
$ws.Range("I136").Value = 1705.9565
$ws.Range("J136").Value = 4171.143
$ws.Range("K136").Value = 5117.8695
$ws.Range("L136").Value = 12513.429
$ws.Range("M136").Value = -2567.8695
$ws.Range("N136").Value = -17613.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3241.75
$ws.Range("I3").Value = 1890.4
$ws.Range("K3").Value = 1890.4
$ws.Range("M3").Value = -1776.4
$ws.Range("H134").Value = 12937.383
$ws.Range("I134").Value = 14993.73
$ws.Range("K134").Value = 44981.19
$ws.Range("M134").Value = -42446.19

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 155.97058
$ws.Range("I7").Value = 91.95999999999999
$ws.Range("J7").Value = 333.77777
$ws.Range("K7").Value = 91.95999999999999
$ws.Range("L7").Value = 333.77777
$ws.Range("M7").Value = 21.04000000000001
$ws.Range("N7").Value = -559.7777699999999
$ws.Range("H22").Value = 874.75
$ws.Range("I22").Value = 833
$ws.Range("K22").Value = 833
$ws.Range("M22").Value = -483
$ws.Range("H31").Value = 4231.727
$ws.Range("I31").Value = 1607.9
$ws.Range("J31").Value = 6418.25
$ws.Range("K31").Value = 1607.9
$ws.Range("L31").Value = 6418.25
$ws.Range("M31").Value = -1312.9
$ws.Range("N31").Value = -7008.25
$ws.Range("H34").Value = 4231.727
$ws.Range("I34").Value = 1607.9
$ws.Range("J34").Value = 6418.25
$ws.Range("K34").Value = 1607.9
$ws.Range("L34").Value = 6418.25
$ws.Range("M34").Value = -1405.9
$ws.Range("N34").Value = -6822.25
$ws.Range("H132").Value = 80362510
$ws.Range("I132").Value = 125004310
$ws.Range("J132").Value = 20840114
$ws.Range("K132").Value = 375012930
$ws.Range("L132").Value = 62520342
$ws.Range("M132").Value = -375010400
$ws.Range("N132").Value = -62525402
$ws.Range("H134").Value = 2428.4092
$ws.Range("I134").Value = 2305.9524
$ws.Range("K134").Value = 6917.8572
$ws.Range("M134").Value = -4382.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3808
$ws.Range("J88").Value = 4000
$ws.Range("L88").Value = 12000
$ws.Range("N88").Value = -12856
$ws.Range("H91").Value = 3808
$ws.Range("J91").Value = 4000
$ws.Range("L91").Value = 12000
$ws.Range("N91").Value = -14964
$ws.Range("H137").Value = 3349.9092
$ws.Range("I137").Value = 1639.8889
$ws.Range("J137").Value = 4533.769
$ws.Range("K137").Value = 4919.6667
$ws.Range("L137").Value = 13601.307
$ws.Range("M137").Value = 180.3333000000002
$ws.Range("N137").Value = -23801.307
$ws.Range("H140").Value = 2017.3448
$ws.Range("I140").Value = 1604.2916
$ws.Range("K140").Value = 4812.8748
$ws.Range("M140").Value = 367.1252000000004
$ws.Range("H141").Value = 2032.3636
$ws.Range("I141").Value = 1875.6
$ws.Range("J141").Value = 3600
$ws.Range("K141").Value = 5626.799999999999
$ws.Range("L141").Value = 10800
$ws.Range("M141").Value = -446.7999999999993
$ws.Range("N141").Value = -21160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 502300
$ws.Range("J80").Value = 145242.14
$ws.Range("L80").Value = 145242.14
$ws.Range("N80").Value = -147238.14
$ws.Range("H83").Value = 502300
$ws.Range("J83").Value = 145242.14
$ws.Range("L83").Value = 726210.7000000001
$ws.Range("N83").Value = -736194.7000000001
$ws.Range("H102").Value = 2674.4546
$ws.Range("I102").Value = 2403.1667
$ws.Range("K102").Value = 2403.1667
$ws.Range("M102").Value = -781.1667000000002
$ws.Range("H132").Value = 18084612
$ws.Range("I132").Value = 27365970
$ws.Range("J132").Value = 10384.895
$ws.Range("K132").Value = 82097910
$ws.Range("L132").Value = 31154.685
$ws.Range("M132").Value = -82095380
$ws.Range("N132").Value = -36214.685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4703.8335
$ws.Range("I16").Value = 7185.533
$ws.Range("J16").Value = 567.6667
$ws.Range("K16").Value = 7185.533
$ws.Range("L16").Value = 567.6667
$ws.Range("M16").Value = -7015.533
$ws.Range("N16").Value = -907.6667
$ws.Range("H132").Value = 1117353.4
$ws.Range("I132").Value = 1282146.8
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 3846440.4
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -3843910.4
$ws.Range("N132").Value = -20055.5
$ws.Range("H136").Value = 4409.6665
$ws.Range("I136").Value = 2194.85
$ws.Range("J136").Value = 7817.077
$ws.Range("K136").Value = 6584.549999999999
$ws.Range("L136").Value = 23451.231
$ws.Range("M136").Value = -4034.549999999999
$ws.Range("N136").Value = -28551.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3699.5
$ws.Range("I62").Value = 3249
$ws.Range("K62").Value = 3249
$ws.Range("M62").Value = -2625
$ws.Range("H65").Value = 3699.5
$ws.Range("I65").Value = 3249
$ws.Range("K65").Value = 16245
$ws.Range("M65").Value = -13125
$ws.Range("H81").Value = 7625.263
$ws.Range("I81").Value = 8212.714
$ws.Range("J81").Value = 7282.5835
$ws.Range("K81").Value = 16425.428
$ws.Range("L81").Value = 14565.167
$ws.Range("M81").Value = -15364.428
$ws.Range("N81").Value = -16687.167
$ws.Range("H84").Value = 7625.263
$ws.Range("I84").Value = 8212.714
$ws.Range("J84").Value = 7282.5835
$ws.Range("K84").Value = 82127.14
$ws.Range("L84").Value = 72825.83499999999
$ws.Range("M84").Value = -76823.14
$ws.Range("N84").Value = -83433.83499999999
$ws.Range("H126").Value = 5496.5713
$ws.Range("I126").Value = 5623.25
$ws.Range("K126").Value = 16869.75
$ws.Range("M126").Value = -14399.75
$ws.Range("H132").Value = 4602299.5
$ws.Range("I132").Value = 6174744
$ws.Range("K132").Value = 18524232
$ws.Range("L132").Value = 17772
$ws.Range("M132").Value = -18521702
$ws.Range("N132").Value = -22832
